$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 57) with the latest metric reading
$ws.Range("A57").Value = "2025-04-29 07:31:56"
$ws.Range("B57").Value = 175
